# Apply the "Add files via upload" commit to user_stories.xlsx
# The edit rewrites several acceptance-criteria cells in the "Sprint 2" sheet
# so that USER STORY 1 ("Manter cadastro de vendas") and USER STORY 2
# ("Confirmar venda") have their own, correct DADO QUE / QUANDO / ENTAO text
# instead of reusing text that actually belongs to the other story.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 2")

# --- USER STORY 1 (A1/A2 = "Manter cadastro de vendas"), CRITERIO DE ACEITE 1 ---
$ws.Range("C5").Value = "DADO QUE: estou logado como vendedor e quero cadastrar uma nova venda"
$ws.Range("C6").Value = "QUANDO: preencho todas as informações obrigatórias para cadastro e insiro um código de produto já existente no sistema"
# C7 ("ENTÃO: o sistema mostra um preview da venda...") is unchanged.

# CRITERIO DE ACEITE 2 (C8/C9/C10) is unchanged.

# --- USER STORY 2 (A13/A14 = "Confirmar venda"), CRITERIO DE ACEITE 1 ---
$ws.Range("C17").Value = 'DADO QUE: realizei o cadastro da venda e seleciono a opção "confirmar"'
$ws.Range("C18").Value = "QUANDO: todos os produtos da compra estão disponíveis no estoque"
# C19 ("ENTÃO: o registro da venda é armazenado no sistema") is unchanged.

# CRITERIO DE ACEITE 2
$ws.Range("C20").Value = 'DADO QUE: realizei o cadastro da venda e seleciono a opção "confirmar"'
$ws.Range("C21").Value = "QUANDO: algum produto da compra não está mais disponível no estoque"
$ws.Range("C22").Value = "ENTÃO: o sistema mostra uma mensagem informando qual produto está em falta e retorna para o cadastro de vendas"

# Recompute the best-fit width of column C now that some text got longer,
# and refresh the view (zoom / scroll / selection) to match the saved file.
$ws.Columns.Item(3).AutoFit()

$ws.Activate()
$excel.ActiveWindow.Zoom = 115
$ws.Range("A12").Select()
$excel.ActiveWindow.ScrollRow = 12
$ws.Range("C22").Select()
